$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "trainingimages/03_kikita"
$ws.Range("B2").Value = "pngimages/03_box.png"
$ws.Range("C2").Value = "trainingimages/11_tokiko"
$ws.Range("D2").Value = "pngimages/11_compass.png"
$ws.Range("E2").Value = 0.5
$ws.Range("F2").Value = -0.5

# Row 3
$ws.Range("A3").Value = "trainingimages/16_kotapi"
$ws.Range("B3").Value = "pngimages/16_icecream.png"
$ws.Range("C3").Value = "trainingimages/10_tokiti"
$ws.Range("D3").Value = "pngimages/10_backpack.png"

# Row 4
$ws.Range("A4").Value = "trainingimages/25_tapapi"
$ws.Range("B4").Value = "pngimages/25_apple.png"
$ws.Range("C4").Value = "trainingimages/18_popata"
$ws.Range("D4").Value = "pngimages/18_donut.png"
$ws.Range("E4").Value = -0.5
$ws.Range("F4").Value = 0.5
